# Commit message: "exclusion to all white milk brands"
# The "Include" sheet's rule rows that filtered on a specific "brand" param
# are renamed to "brand_name" so the rule now applies to all (white milk)
# brands rather than a single one. The workbook is also left with the
# "Include" tab active/selected, matching the saved UI state.

$wb = $excel.ActiveWorkbook

$includeSheet = $wb.Worksheets.Item("Include")

# Rename "brand" -> "brand_name" wherever it appears as a Param 1 value
# on the Include sheet (rows 2 and 4, column C).
$usedRange = $includeSheet.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $includeSheet.Cells.Item($r, $c)
        if ($cell.Value2 -eq "brand") {
            $cell.Value = "brand_name"
        }
    }
}

# Make the "Include" sheet the active tab (was "Exclude").
$includeSheet.Activate()
